# Updated cryptos list with refreshed prices and volumes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.376.20'
$ws.Range('E2').Value = '  -4.73%  '
$ws.Range('D3').Value = '3.699.91'
$ws.Range('E3').Value = '  -4.81%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.07%  '
$ws.Range('D7').Value = '3.694.30'
$ws.Range('E7').Value = '  -4.83%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.627'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  -5.55%  '
$ws.Range('E11').Value = '  -8.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.88'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000291'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -9.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.89%  '
$ws.Range('D15').Value = '4.193.70'
$ws.Range('E15').Value = '  -7.32%  '
$ws.Range('D16').Value = '3.704.52'
$ws.Range('E16').Value = '  -4.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.29%  '
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('E19').Value = '  -7.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.45%  '
$ws.Range('D21').Value = '68.199.44'
$ws.Range('E21').Value = '  -4.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '408.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.68%  '
$ws.Range('E25').Value = '  -8.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.89%  '
$ws.Range('E28').Value = '  -4.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.06'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.95%  '
$ws.Range('E30').Value = '  -7.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.64'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.00%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '65.51'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.07%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.117'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '43.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -15.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '597.94'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.97%  '
$ws.Range('D38').Value = '0.0₃0889'
$ws.Range('E38').Value = '  -10.02%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('E40').Value = '  -4.63%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -4.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.88%  '
$ws.Range('E44').Value = '  -9.13%  '
$ws.Range('E45').Value = '  -8.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0435'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.53%  '
$ws.Range('E47').Value = '  -9.57%  '
$ws.Range('D48').Value = '2.784.56'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.17'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.46%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.133'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.61%  '
$ws.Range('E51').Value = '  -3.49%  '
